$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new entry "Richard Zhao" at G5, pushing the existing G5:G12
# values down to G6:G13 (column-only shift; columns A-F and H are
# untouched, matching the source diff). Read with Value2 (plain .Value
# reads misbehave in this host) and write bottom-up so no source cell
# is clobbered before it has been read.
$g12 = $ws.Range("G12").Value2
$g11 = $ws.Range("G11").Value2
$g10 = $ws.Range("G10").Value2
$g9  = $ws.Range("G9").Value2
$g8  = $ws.Range("G8").Value2
$g7  = $ws.Range("G7").Value2
$g6  = $ws.Range("G6").Value2
$g5  = $ws.Range("G5").Value2

$ws.Range("G13").Value = $g12
$ws.Range("G12").Value = $g11
$ws.Range("G11").Value = $g10
$ws.Range("G10").Value = $g9
$ws.Range("G9").Value  = $g8
$ws.Range("G8").Value  = $g7
$ws.Range("G7").Value  = $g6
$ws.Range("G6").Value  = $g5
$ws.Range("G5").Value  = "Richard Zhao"

# Update the active selection on the sheet view to F4.
$ws.Range("F4").Select() | Out-Null
